$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Delete rows 5 and 6 (matches from 6 rows of data to 4)
$ws.Rows("5:6").Delete()

# Row 2: update match data
$ws.Range("B2").Value = '2025-11-29T13:30:00'
$ws.Range("C2").Value = 'Сибирь'
$ws.Range("D2").Value = 'Автомобилист'
$ws.Range("E2").Value = 897828
$ws.Range("F2").Value = 'https://text.khl.ru/text/897828.html'
$ws.Range("G2").Value = 0.833333
$ws.Range("H2").Value = 3.061201
$ws.Range("I2").Value = 2.620406
$ws.Range("J2").Value = 2.68588
$ws.Range("K2").Value = 1.759607
$ws.Range("L2").Value = 2.840804
$ws.Range("M2").Value = 3.894535
$ws.Range("N2").Value = 22.523066
$ws.Range("O2").Value = 33.006542
$ws.Range("P2").Value = 55.529608
$ws.Range("Q2").Value = -0.2
$ws.Range("R2").Value = 0.05918
$ws.Range("S2").Value = 0.225401
$ws.Range("T2").Value = 0.171304
$ws.Range("U2").Value = 0.603107
$ws.Range("V2").Value = 0.325639
$ws.Range("W2").Value = 0.674173
$ws.Range("X2").Value = 0.513157
$ws.Range("Y2").Value = 0.486655
$ws.Range("Z2").Value = 0.685689
$ws.Range("AA2").Value = 0.314124
$ws.Range("AB2").Value = 0.817975
$ws.Range("AC2").Value = 0.181838
$ws.Range("AD2").Value = 0.904913
$ws.Range("AE2").Value = 0.094899
$ws.Range("AF2").Value = 0.525037
$ws.Range("AG2").Value = 0.474963
$ws.Range("AH2").Value = 0.258588
$ws.Range("AI2").Value = 0.741412
$ws.Range("AJ2").Value = 0.775779
$ws.Range("AK2").Value = 0.224221
$ws.Range("AL2").Value = 0.540216
$ws.Range("AM2").Value = 0.459784
$ws.Range("AN2").Value = 0.588135
$ws.Range("AO2").Value = 0.892984

# Row 3: update match data
$ws.Range("B3").Value = '2025-11-29T14:30:00'
$ws.Range("C3").Value = 'Металлург Мг'
$ws.Range("D3").Value = 'Барыс'
$ws.Range("E3").Value = 897829
$ws.Range("F3").Value = 'https://text.khl.ru/text/897829.html'
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 1.808836
$ws.Range("I3").Value = 4.59375
$ws.Range("J3").Value = 4.736323
$ws.Range("K3").Value = 5.368161
$ws.Range("L3").Value = 3.201293
$ws.Range("M3").Value = 7.808836
$ws.Range("N3").Value = 38.880765
$ws.Range("O3").Value = 26.496915
$ws.Range("P3").Value = 65.37768
$ws.Range("Q3").Value = 0.2
$ws.Range("R3").Value = -0.066061
$ws.Range("S3").Value = 0.693818
$ws.Range("T3").Value = 0.106573
$ws.Range("U3").Value = 0.177431
$ws.Range("V3").Value = 0.028695
$ws.Range("W3").Value = 0.949128
$ws.Range("X3").Value = 0.071346
$ws.Range("Y3").Value = 0.906476
$ws.Range("Z3").Value = 0.144446
$ws.Range("AA3").Value = 0.833376
$ws.Range("AB3").Value = 0.248851
$ws.Range("AC3").Value = 0.728972
$ws.Range("AD3").Value = 0.376664
$ws.Range("AE3").Value = 0.601159
$ws.Range("AF3").Value = 0.970307
$ws.Range("AG3").Value = 0.029693
$ws.Range("AH3").Value = 0.903124
$ws.Range("AI3").Value = 0.096876
$ws.Range("AJ3").Value = 0.828967
$ws.Range("AK3").Value = 0.171033
$ws.Range("AL3").Value = 0.620366
$ws.Range("AM3").Value = 0.379634
$ws.Range("AN3").Value = 0.877549
$ws.Range("AO3").Value = 0.413387

# Row 4: update match data
$ws.Range("B4").Value = '2025-11-29T17:00:00'
$ws.Range("C4").Value = 'Спартак'
$ws.Range("D4").Value = 'Лада'
$ws.Range("E4").Value = 897830
$ws.Range("F4").Value = 'https://text.khl.ru/text/897830.html'
$ws.Range("G4").Value = 4.421053
$ws.Range("H4").Value = 1.111111
$ws.Range("I4").Value = 3.683199
$ws.Range("J4").Value = 2.541289
$ws.Range("K4").Value = 3.481171
$ws.Range("L4").Value = 2.397155
$ws.Range("M4").Value = 5.532164
$ws.Range("N4").Value = 41.634955
$ws.Range("O4").Value = 21.960244
$ws.Range("P4").Value = 63.595199
$ws.Range("Q4").Value = 0.2
$ws.Range("R4").Value = -0.2
$ws.Range("S4").Value = 0.591843
$ws.Range("T4").Value = 0.153703
$ws.Range("U4").Value = 0.253435
$ws.Range("V4").Value = 0.162396
$ws.Range("W4").Value = 0.836584
$ws.Range("X4").Value = 0.301673
$ws.Range("Y4").Value = 0.697308
$ws.Range("Z4").Value = 0.465416
$ws.Range("AA4").Value = 0.533565
$ws.Range("AB4").Value = 0.625838
$ws.Range("AC4").Value = 0.373142
$ws.Range("AD4").Value = 0.760555
$ws.Range("AE4").Value = 0.238426
$ws.Range("AF4").Value = 0.862108
$ws.Range("AG4").Value = 0.137892
$ws.Range("AH4").Value = 0.675656
$ws.Range("AI4").Value = 0.324344
$ws.Range("AJ4").Value = 0.690939
$ws.Range("AK4").Value = 0.309061
$ws.Range("AL4").Value = 0.429548
$ws.Range("AM4").Value = 0.570452
$ws.Range("AN4").Value = 0.861453
$ws.Range("AO4").Value = 0.575461

